$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema3e"
$ws.Cells.Item(2,3).Value = "Plxnd1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.03046233333333333
$ws.Cells.Item(2,8).Value = 0.091387
$ws.Cells.Item(2,9).Value = 0.01378316480238178
$ws.Cells.Item(2,10).Value = 0.01378316480238178
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 88.37814633333333
$ws.Cells.Item(2,14).Value = 265.134439
$ws.Cells.Item(2,15).Value = 0.7138016014383547
$ws.Cells.Item(2,16).Value = 0.7138016014383547
$ws.Cells.Item(2,17).Value = 2.692204552988111
$ws.Cells.Item(2,18).Value = 24.229840976893
$ws.Cells.Item(2,19).Value = 0.009838445108828878
$ws.Cells.Item(2,20).Value = 0.009838445108828876

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema3e"
$ws.Cells.Item(3,3).Value = "Plxnd1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.03046233333333333
$ws.Cells.Item(3,8).Value = 0.091387
$ws.Cells.Item(3,9).Value = 0.01378316480238178
$ws.Cells.Item(3,10).Value = 0.01378316480238178
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 8.866675333333333
$ws.Cells.Item(3,14).Value = 26.600026
$ws.Cells.Item(3,15).Value = 0.07161325864989523
$ws.Cells.Item(3,16).Value = 0.07161325864989525
$ws.Cells.Item(3,17).Value = 0.2700996195624444
$ws.Cells.Item(3,18).Value = 2.430896576062
$ws.Cells.Item(3,19).Value = 0.0009870573460070986
$ws.Cells.Item(3,20).Value = 0.0009870573460070986

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema3e"
$ws.Cells.Item(4,3).Value = "Plxnd1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.03046233333333333
$ws.Cells.Item(4,8).Value = 0.091387
$ws.Cells.Item(4,9).Value = 0.01378316480238178
$ws.Cells.Item(4,10).Value = 0.01378316480238178
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 26.56849866666667
$ws.Cells.Item(4,14).Value = 79.705496
$ws.Cells.Item(4,15).Value = 0.2145851399117501
$ws.Cells.Item(4,16).Value = 0.2145851399117501
$ws.Cells.Item(4,17).Value = 0.8093384625502221
$ws.Cells.Item(4,18).Value = 7.284046162951999
$ws.Cells.Item(4,19).Value = 0.002957662347545804
$ws.Cells.Item(4,20).Value = 0.002957662347545803

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Sema3e"
$ws.Cells.Item(5,3).Value = "Plxnd1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.6579096666666667
$ws.Cells.Item(5,8).Value = 1.973729
$ws.Cells.Item(5,9).Value = 0.2976816405204262
$ws.Cells.Item(5,10).Value = 0.2976816405204262
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 88.37814633333333
$ws.Cells.Item(5,14).Value = 265.134439
$ws.Cells.Item(5,15).Value = 0.7138016014383547
$ws.Cells.Item(5,16).Value = 0.7138016014383547
$ws.Cells.Item(5,17).Value = 58.14483679478123
$ws.Cells.Item(5,18).Value = 523.303531153031
$ws.Cells.Item(5,19).Value = 0.2124856317222769
$ws.Cells.Item(5,20).Value = 0.2124856317222768

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema3e"
$ws.Cells.Item(6,3).Value = "Plxnd1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.6579096666666667
$ws.Cells.Item(6,8).Value = 1.973729
$ws.Cells.Item(6,9).Value = 0.2976816405204262
$ws.Cells.Item(6,10).Value = 0.2976816405204262
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 8.866675333333333
$ws.Cells.Item(6,14).Value = 26.600026
$ws.Cells.Item(6,15).Value = 0.07161325864989523
$ws.Cells.Item(6,16).Value = 0.07161325864989525
$ws.Cells.Item(6,17).Value = 5.833471412994889
$ws.Cells.Item(6,18).Value = 52.501242716954
$ws.Cells.Item(6,19).Value = 0.02131795231791442
$ws.Cells.Item(6,20).Value = 0.02131795231791442

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema3e"
$ws.Cells.Item(7,3).Value = "Plxnd1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.6579096666666667
$ws.Cells.Item(7,8).Value = 1.973729
$ws.Cells.Item(7,9).Value = 0.2976816405204262
$ws.Cells.Item(7,10).Value = 0.2976816405204262
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 26.56849866666667
$ws.Cells.Item(7,14).Value = 79.705496
$ws.Cells.Item(7,15).Value = 0.2145851399117501
$ws.Cells.Item(7,16).Value = 0.2145851399117501
$ws.Cells.Item(7,17).Value = 17.47967210162044
$ws.Cells.Item(7,18).Value = 157.317048914584
$ws.Cells.Item(7,19).Value = 0.06387805648023497
$ws.Cells.Item(7,20).Value = 0.06387805648023495

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Sema3e"
$ws.Cells.Item(8,3).Value = "Plxnd1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.521739666666667
$ws.Cells.Item(8,8).Value = 4.565219
$ws.Cells.Item(8,9).Value = 0.688535194677192
$ws.Cells.Item(8,10).Value = 0.6885351946771919
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 88.37814633333333
$ws.Cells.Item(8,14).Value = 265.134439
$ws.Cells.Item(8,15).Value = 0.7138016014383547
$ws.Cells.Item(8,16).Value = 0.7138016014383547
$ws.Cells.Item(8,17).Value = 134.4885309419046
$ws.Cells.Item(8,18).Value = 1210.396778477141
$ws.Cells.Item(8,19).Value = 0.491477524607249
$ws.Cells.Item(8,20).Value = 0.4914775246072489

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Sema3e"
$ws.Cells.Item(9,3).Value = "Plxnd1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.521739666666667
$ws.Cells.Item(9,8).Value = 4.565219
$ws.Cells.Item(9,9).Value = 0.688535194677192
$ws.Cells.Item(9,10).Value = 0.6885351946771919
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 8.866675333333333
$ws.Cells.Item(9,14).Value = 26.600026
$ws.Cells.Item(9,15).Value = 0.07161325864989523
$ws.Cells.Item(9,16).Value = 0.07161325864989525
$ws.Cells.Item(9,17).Value = 13.49277156618822
$ws.Cells.Item(9,18).Value = 121.434944095694
$ws.Cells.Item(9,19).Value = 0.04930824898597373
$ws.Cells.Item(9,20).Value = 0.04930824898597373

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sema3e"
$ws.Cells.Item(10,3).Value = "Plxnd1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.521739666666667
$ws.Cells.Item(10,8).Value = 4.565219
$ws.Cells.Item(10,9).Value = 0.688535194677192
$ws.Cells.Item(10,10).Value = 0.6885351946771919
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 26.56849866666667
$ws.Cells.Item(10,14).Value = 79.705496
$ws.Cells.Item(10,15).Value = 0.2145851399117501
$ws.Cells.Item(10,16).Value = 0.2145851399117501
$ws.Cells.Item(10,17).Value = 40.43033830484711
$ws.Cells.Item(10,18).Value = 363.873044743624
$ws.Cells.Item(10,19).Value = 0.1477494210839694
$ws.Cells.Item(10,20).Value = 0.1477494210839693
